# Refresh the cryptocurrency price/volume figures (columns D and E)
# for each changed row, mirroring the source site's latest snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ Price = <new D value or $null>; Volume = <new E value or $null> }
$updates = @{
    2 = @{ Price = "62.890.46"; Volume = "  +2.02%  " }
    3 = @{ Price = "3.474.09"; Volume = "  +2.33%  " }
    4 = @{ Price = "0.998"; Volume = "  -0.12%  " }
    5 = @{ Price = "579.50"; Volume = "  +0.53%  " }
    6 = @{ Price = "147.99"; Volume = "  +4.26%  " }
    7 = @{ Price = $null; Volume = "  -0.13%  " }
    8 = @{ Price = "0.482"; Volume = "  +1.75%  " }
    9 = @{ Price = "7.61"; Volume = "  -0.99%  " }
    10 = @{ Price = "0.125"; Volume = "  +2.00%  " }
    11 = @{ Price = "0.402"; Volume = "  +4.28%  " }
    12 = @{ Price = "4.070.39"; Volume = "  +2.45%  " }
    13 = @{ Price = "29.98"; Volume = "  +5.89%  " }
    14 = @{ Price = "0.128"; Volume = "  +2.48%  " }
    15 = @{ Price = "3.474.88"; Volume = "  +2.61%  " }
    16 = @{ Price = $null; Volume = "  +0.84%  " }
    17 = @{ Price = "62.831.11"; Volume = "  +1.96%  " }
    18 = @{ Price = "6.36"; Volume = "  +3.49%  " }
    19 = @{ Price = "14.41"; Volume = "  +5.72%  " }
    20 = @{ Price = "9.29"; Volume = "  +3.03%  " }
    21 = @{ Price = "388.83"; Volume = "  -0.58%  " }
    22 = @{ Price = "0.562"; Volume = "  +2.28%  " }
    23 = @{ Price = "74.61"; Volume = "  -0.25%  " }
    24 = @{ Price = $null; Volume = "  +0.03%  " }
    25 = @{ Price = "3.610.82"; Volume = "  +2.34%  " }
    26 = @{ Price = $null; Volume = "  +2.22%  " }
    27 = @{ Price = $null; Volume = "  -9.70%  " }
    28 = @{ Price = "7.60"; Volume = "  +2.78%  " }
    29 = @{ Price = $null; Volume = "  -0.07%  " }
    30 = @{ Price = "8.20"; Volume = "  +2.48%  " }
    31 = @{ Price = $null; Volume = "  -0.17%  " }
    32 = @{ Price = $null; Volume = "  +0.06%  " }
    33 = @{ Price = "1.40"; Volume = "  -1.65%  " }
    34 = @{ Price = "23.77"; Volume = "  +1.92%  " }
    35 = @{ Price = "5.32"; Volume = "  +5.47%  " }
    36 = @{ Price = "7.12"; Volume = "  +2.83%  " }
    37 = @{ Price = "1.60"; Volume = "  +8.12%  " }
    38 = @{ Price = "31.55"; Volume = "  +22.23%  " }
    39 = @{ Price = "169.47"; Volume = "  +0.17%  " }
    40 = @{ Price = "3.511.66"; Volume = "  +2.42%  " }
    41 = @{ Price = "0.0766"; Volume = "  -0.04%  " }
    42 = @{ Price = "0.800"; Volume = "  +2.54%  " }
    43 = @{ Price = "42.41"; Volume = "  -0.06%  " }
    44 = @{ Price = "4.49"; Volume = "  +1.33%  " }
    45 = @{ Price = "1.73"; Volume = "  +4.82%  " }
    46 = @{ Price = $null; Volume = "  +3.72%  " }
    47 = @{ Price = "2.614.80"; Volume = "  +6.28%  " }
    48 = @{ Price = "23.21"; Volume = "  +1.84%  " }
    49 = @{ Price = "2.26"; Volume = "  +12.09%  " }
    50 = @{ Price = "6.76"; Volume = "  +1.48%  " }
    51 = @{ Price = "0.998"; Volume = "  -0.10%  " }
}

foreach ($row in $updates.Keys) {
    $entry = $updates[$row]

    if ($null -ne $entry.Price) {
        $priceCell = $ws.Range("D$row")
        # Prices like "579.50" / "0.998" parse as numbers unless the cell is
        # forced to Text first; switch back to the default "Normal" style
        # afterwards so no stray number-format style sticks to the cell
        # (matches values such as "62.890.46" that already stay textual).
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $entry.Price
        $priceCell.Style = "Normal"
    }

    if ($null -ne $entry.Volume) {
        $ws.Range("E$row").Value = $entry.Volume
    }
}
